$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("body")

# Update the xpath selector text used by the "contact_button" rows (3, 5, 7, 9)
# from the old absolute xpath to the new relative xpath that targets the
# showcase image.
$oldXpath = "/html/body/div[1]/div/div/div[5]/div"
$newXpath = '//*[@id="root"]/div/div/div[4]/div/div[1]/img'

$usedRange = $ws.UsedRange
for ($r = 1; $r -le $usedRange.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldXpath) {
        $cell.Value = $newXpath
    }
}

# Update the selected cell in the sheet view to C9
$ws.Range("C9").Select()
